$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1001924303428215
$ws.Range("C2").Value = 0.1001924303428215
$ws.Range("D2").Value = 0.01289511287924234
$ws.Range("F2").Value = 0.01074592739936861
$ws.Range("H2").Value = 0.000002072645637406635
$ws.Range("I2").Value = 0.00000004333390484538357
$ws.Range("J2").Value = 0.00214706950033147
$ws.Range("K2").Value = 0.00000414529127481327
$ws.Range("L2").Value = 0.0168702991656713
$ws.Range("N2").Value = 0.2430495309010738
$ws.Range("B3").Value = 0.1001924303428215
$ws.Range("C3").Value = 0.09418088452225222
$ws.Range("D3").Value = 0.01697239176792528
$ws.Range("E3").Value = 0.0000154741354550908
$ws.Range("F3").Value = 0.001531939410053989
$ws.Range("G3").Value = 0.0000005802800795659052
$ws.Range("H3").Value = 0.00005748793720239029
$ws.Range("I3").Value = 0.000003330374339716509
$ws.Range("J3").Value = 0.0003880168470961185
$ws.Range("K3").Value = 0.00000414646822366887
$ws.Range("L3").Value = 0.00005874844645855104
$ws.Range("N3").Value = 0.2134054305319081
$ws.Range("B4").Value = 0.1001924303428215
$ws.Range("C4").Value = 0.08853003145091709
$ws.Range("D4").Value = 0.02007738607407162
$ws.Range("E4").Value = 0.0000354077297838586
$ws.Range("F4").Value = 0.002016770439371266
$ws.Range("G4").Value = 0.000001327789866894697
$ws.Range("H4").Value = 0.00007672727130783347
$ws.Range("I4").Value = 0.000007565083759025392
$ws.Range("J4").Value = 0.0004397593204317707
$ws.Range("K4").Value = 0.000004673121595079125
$ws.Range("L4").Value = 0.00007523043822024778
$ws.Range("N4").Value = 0.2114573090621462
$ws.Range("B5").Value = 0.1001924303428215
$ws.Range("C5").Value = 0.08321822956386206
$ws.Range("D5").Value = 0.0223775800500159
$ws.Range("E5").Value = 0.0000585091766387965
$ws.Range("F5").Value = 0.002386223830536418
$ws.Range("G5").Value = 0.000002194094123954869
$ws.Range("H5").Value = 0.00009093326402336103
$ws.Range("I5").Value = 0.00001247299486266094
$ws.Range("J5").Value = 0.0005196506076645026
$ws.Range("K5").Value = 0.000005876697764851445
$ws.Range("L5").Value = 0.00008899489842649483
$ws.Range("N5").Value = 0.2089530955207405
$ws.Range("B6").Value = 0.1001924303428215
$ws.Range("C6").Value = 0.07822513579003033
$ws.Range("D6").Value = 0.02401403681634524
$ws.Range("E6").Value = 0.00008372401575292929
$ws.Range("F6").Value = 0.002660159126982078
$ws.Range("G6").Value = 0.000003139650590734848
$ws.Range("H6").Value = 0.000101456224444874
$ws.Range("I6").Value = 0.00001783003391926049
$ws.Range("J6").Value = 0.0005800438410456883
$ws.Range("K6").Value = 0.000007864779485592351
$ws.Range("L6").Value = 0.00009923521352973044
$ws.Range("N6").Value = 0.2059850558349479
$ws.Range("B7").Value = 0.1001924303428215
$ws.Range("C7").Value = 0.07353162764262851
$ws.Range("D7").Value = 0.02510543944129527
$ws.Range("E7").Value = 0.0001101965874914616
$ws.Range("F7").Value = 0.002855303942640225
$ws.Range("G7").Value = 0.000004132372030929808
$ws.Range("H7").Value = 0.0001089507847710486
$ws.Range("I7").Value = 0.00002345437235768923
$ws.Range("J7").Value = 0.000623113635203065
$ws.Range("K7").Value = 0.00001070913948689879
$ws.Range("L7").Value = 0.0001065315430307197
$ws.Range("N7").Value = 0.2026718898037573
$ws.Range("B8").Value = 0.1001924303428215
$ws.Range("C8").Value = 0.0691197299840708
$ws.Range("D8").Value = 0.02575152118365869
$ws.Range("E8").Value = 0.000137237610371255
$ws.Range("F8").Value = 0.00298573292632188
$ws.Range("G8").Value = 0.00000514641038892206
$ws.Range("H8").Value = 0.0001139583779355288
$ws.Range("I8").Value = 0.00002919953698920423
$ws.Range("J8").Value = 0.0006519194254799083
$ws.Range("K8").Value = 0.00001445249184493044
$ws.Range("L8").Value = 0.0001114087381127642
$ws.Range("N8").Value = 0.1991127370279954
$ws.Range("B9").Value = 0.1001924303428215
$ws.Range("C9").Value = 0.06497254618502656
$ws.Range("D9").Value = 0.02603597680515414
$ws.Range("E9").Value = 0.0001642967827012502
$ws.Range("F9").Value = 0.003063274331080456
$ws.Range("G9").Value = 0.000006161129351296884
$ws.Range("H9").Value = 0.0001169338194051487
$ws.Range("I9").Value = 0.00003494858785533948
$ws.Range("J9").Value = 0.0006690680075629291
$ws.Range("K9").Value = 0.00001911414226831253
$ws.Range("L9").Value = 0.0001143089465474045
$ws.Range("N9").Value = 0.1953890590797743
$ws.Range("B10").Value = 0.1001924303428215
$ws.Range("C10").Value = 0.06107419341392496
$ws.Range("D10").Value = 0.02602893305548261
$ws.Range("E10").Value = 0.0001909396449518002
$ws.Range("F10").Value = 0.003097855080828918
$ws.Range("G10").Value = 0.000007160236685692507
$ws.Range("H10").Value = 0.0001182587927634216
$ws.Range("I10").Value = 0.00004060920148070428
$ws.Range("J10").Value = 0.0006767481134661974
$ws.Range("K10").Value = 0.00002469470792202909
$ws.Range("L10").Value = 0.0001156032716841219
$ws.Range("N10").Value = 0.1915674258620119
$ws.Range("B11").Value = 0.1001924303428215
$ws.Range("C11").Value = 0.05740974180908946
$ws.Range("D11").Value = 0.02578904470199572
$ws.Range("E11").Value = 0.0002168280545597289
$ws.Range("F11").Value = 0.003097793590659432
$ws.Range("G11").Value = 0.000008131052045989835
$ws.Range("H11").Value = 0.0001182531173417854
$ws.Range("I11").Value = 0.00004610952176839262
$ws.Range("J11").Value = 0.0006767937850875247
$ws.Range("K11").Value = 0.00003118005148579132
$ws.Range("L11").Value = 0.0001156026836192444
$ws.Range("N11").Value = 0.1877019087104746
$ws.Range("B12").Value = 0.1001924303428215
$ws.Range("C12").Value = 0.0539651573005441
$ws.Range("D12").Value = 0.02536527250524172
$ws.Range("E12").Value = 0.0002417037226744514
$ws.Range("F12").Value = 0.003070048206984779
$ws.Range("G12").Value = 0.000009063889600291926
$ws.Range("H12").Value = 0.0001171843046238216
$ws.Range("I12").Value = 0.00005139466155788653
$ws.Range("J12").Value = 0.0006707391326479288
$ws.Range("K12").Value = 0.00003854455212700506
$ws.Range("L12").Value = 0.0001145673062814577
$ws.Range("N12").Value = 0.1838361059251049
$ws.Range("B13").Value = 0.1001924303428215
$ws.Range("C13").Value = 0.05072724786251145
$ws.Range("D13").Value = 0.02479839106748811
$ws.Range("E13").Value = 0.0002653743454458568
$ws.Range("F13").Value = 0.003020427951952543
$ws.Range("G13").Value = 0.000009951537954219629
$ws.Range("H13").Value = 0.000115275664520554
$ws.Range("I13").Value = 0.00005642375552310235
$ws.Range("J13").Value = 0.0006598648131548821
$ws.Range("K13").Value = 0.00004675381737338737
$ws.Range("L13").Value = 0.0001127142962436338
$ws.Range("N13").Value = 0.1800048554549892
$ws.Range("B14").Value = 0.1001924303428215
$ws.Range("C14").Value = 0.04768361299076077
$ws.Range("D14").Value = 0.02412226727911558
$ws.Range("E14").Value = 0.0002877019330543586
$ws.Range("F14").Value = 0.002953771252270062
$ws.Range("G14").Value = 0.00001078882248953844
$ws.Range("H14").Value = 0.0001127131800973583
$ws.Range("I14").Value = 0.00006116748009066496
$ws.Range("J14").Value = 0.0006452374485615445
$ws.Range("K14").Value = 0.00005576692399716809
$ws.Range("L14").Value = 0.0001102245246268954
$ws.Range("N14").Value = 0.1762356821778855
$ws.Range("B15").Value = 0.1001924303428215
$ws.Range("C15").Value = 0.04482259621131512
$ws.Range("D15").Value = 0.02336494396669389
$ws.Range("E15").Value = 0.0003085929996637752
$ws.Range("F15").Value = 0.002874097479010813
$ws.Range("G15").Value = 0.00001157223748739157
$ws.Range("H15").Value = 0.0001096513364223486
$ws.Range("I15").Value = 0.00006560596880430081
$ws.Range("J15").Value = 0.0006277430472144137
$ws.Range("K15").Value = 0.00006553826254729987
$ws.Range("L15").Value = 0.0001072482425152038
$ws.Range("N15").Value = 0.172550020094496
$ws.Range("B16").Value = 0.1001924303428215
$ws.Range("C16").Value = 0.04213324043863621
$ws.Range("D16").Value = 0.02254955814436871
$ws.Range("E16").Value = 0.0003279903284332221
$ws.Range("F16").Value = 0.00278473539953345
$ws.Range("G16").Value = 0.00001229963731624583
$ws.Range("H16").Value = 0.0001062180612982148
$ws.Range("I16").Value = 0.00006972706238928733
$ws.Range("J16").Value = 0.0006081153331131266
$ws.Range("K16").Value = 0.000007601904872754684
$ws.Range("L16").Value = 0.0001039098832053616
$ws.Range("N16").Value = 0.1689642436798429
$ws.Range("B17").Value = 0.03658806804351665
$ws.Range("C17").Value = 0.1032096083116229
$ws.Range("D17").Value = 0.02169511884903157
$ws.Range("E17").Value = 0.0003458660690103344
$ws.Range("F17").Value = 0.002688432026108409
$ws.Range("G17").Value = 0.00001296997758788754
$ws.Range("H17").Value = 0.000102518911922348
$ws.Range("I17").Value = 0.00007352484196920131
$ws.Range("J17").Value = 0.0005869597515481801
$ws.Range("K17").Value = 0.00008715855511024398
$ws.Range("L17").Value = 0.0001003121315681446
$ws.Range("N17").Value = 0.1654905374689959
$ws.Range("B18").Value = 0.01418469763434395
$ws.Range("C18").Value = 0.1194204022220982
$ws.Range("D18").Value = 0.02463342752037421
$ws.Range("E18").Value = 0.0003622159616968829
$ws.Range("F18").Value = 0.002587444821873152
$ws.Range("G18").Value = 0.00001358309856363311
$ws.Range("H18").Value = 0.0000986406213637292
$ws.Range("I18").Value = 0.00007699840170056671
$ws.Range("J18").Value = 0.000564773804113854
$ws.Range("K18").Value = 0.00009890510843745981
$ws.Range("L18").Value = 0.00009653937120688473
$ws.Range("N18").Value = 0.1621376285657725
$ws.Range("B19").Value = 0.005222413989233897
$ws.Range("C19").Value = 0.1212174617338824
$ws.Range("D19").Value = 0.02810363752564397
$ws.Range("E19").Value = 0.0003816340277938192
$ws.Range("F19").Value = 0.002936991673905836
$ws.Range("G19").Value = 0.00001431127604226822
$ws.Range("H19").Value = 0.0001116555476577235
$ws.Range("I19").Value = 0.00008112397047064627
$ws.Range("J19").Value = 0.0006383057954515601
$ws.Range("K19").Value = 0.0001112068907719206
$ws.Range("L19").Value = 0.0001094920350338763
$ws.Range("N19").Value = 0.1589282344658879
$ws.Range("B20").Value = 0.001777182209405854
$ws.Range("C20").Value = 0.1173896458096775
$ws.Range("D20").Value = 0.03116113960083032
$ws.Range("E20").Value = 0.000404672640046365
$ws.Range("F20").Value = 0.003349817688255304
$ws.Range("G20").Value = 0.00001517522400173869
$ws.Range("H20").Value = 0.0001273678228013037
$ws.Range("I20").Value = 0.00008601892686650739
$ws.Range("J20").Value = 0.0007280510279825225
$ws.Range("K20").Value = 0.0001241680994193833
$ws.Range("L20").Value = 0.0001248834179431609
$ws.Range("N20").Value = 0.15528812246723
$ws.Range("B21").Value = 0.0005586466071316881
$ws.Range("C21").Value = 0.111564802663371
$ws.Range("D21").Value = 0.03353034740928642
$ws.Range("E21").Value = 0.0004307351736460632
$ws.Range("F21").Value = 0.003713719358403991
$ws.Range("G21").Value = 0.00001615256901172737
$ws.Range("H21").Value = 0.0001412938126840444
$ws.Range("I21").Value = 0.00009155648751245692
$ws.Range("J21").Value = 0.0008078492541372528
$ws.Range("K21").Value = 0.0001379116905233447
$ws.Range("L21").Value = 0.0001384729346299182
$ws.Range("N21").Value = 0.1511314879603379
$ws.Range("B22").Value = 0.0001640771271126538
$ws.Range("C22").Value = 0.1052654839835878
$ws.Range("D22").Value = 0.03519468345769572
$ws.Range("E22").Value = 0.0004589110056751172
$ws.Range("F22").Value = 0.003995939665776327
$ws.Range("G22").Value = 0.00001720916271281689
$ws.Range("H22").Value = 0.0001521187976744415
$ws.Range("I22").Value = 0.00009754314445147617
$ws.Range("J22").Value = 0.0008699820336428298
$ws.Range("K22").Value = 0.0001525403585340919
$ws.Range("L22").Value = 0.0001490200269360635
$ws.Range("N22").Value = 0.1465175087637993
